$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.8775636666666666
$ws.Range("H2").Value = 2.632691
$ws.Range("I2").Value = 0.1887436506618166
$ws.Range("J2").Value = 0.2083714858314108
$ws.Range("M2").Value = 1.143899333333333
$ws.Range("N2").Value = 3.431698
$ws.Range("O2").Value = 0.0786649926114418
$ws.Range("P2").Value = 0.09278552072412605
$ws.Range("Q2").Value = 1.003844493257556
$ws.Range("R2").Value = 9.034600439318
$ws.Range("S2").Value = 0.01484751788476836
$ws.Range("T2").Value = 0.0193338568169273
$ws.Range("G3").Value = 0.8775636666666666
$ws.Range("H3").Value = 2.632691
$ws.Range("I3").Value = 0.1887436506618166
$ws.Range("J3").Value = 0.2083714858314108
$ws.Range("O3").Value = 0.4082467010060103
$ws.Range("P3").Value = 0.4815278242489749
$ws.Range("Q3").Value = 5.209638863372112
$ws.Range("R3").Value = 46.886749770349
$ws.Range("S3").Value = 0.07705397271851749
$ws.Range("T3").Value = 0.1003366682079253
$ws.Range("G4").Value = 0.8775636666666666
$ws.Range("H4").Value = 2.632691
$ws.Range("I4").Value = 0.1887436506618166
$ws.Range("J4").Value = 0.2083714858314108
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02736833333333333
$ws.Range("N4").Value = 0.082105
$ws.Range("O4").Value = 0.001882097206211744
$ws.Range("P4").Value = 0.002219937529192361
$ws.Range("Q4").Value = 0.02401745495055555
$ws.Range("R4").Value = 0.216157094555
$ws.Range("S4").Value = 0.0003552338976008103
$ws.Range("T4").Value = 0.0004625716814107231
$ws.Range("G5").Value = 0.8775636666666666
$ws.Range("H5").Value = 2.632691
$ws.Range("I5").Value = 0.1887436506618166
$ws.Range("J5").Value = 0.2083714858314108
$ws.Range("M5").Value = 6.638933
$ws.Range("N5").Value = 13.277866
$ws.Range("O5").Value = 0.4565538244270245
$ws.Range("P5").Value = 0.3590041171790666
$ws.Range("Q5").Value = 5.826086386234333
$ws.Range("R5").Value = 34.95651831740599
$ws.Range("S5").Value = 0.08617163554597065
$ws.Range("T5").Value = 0.07480622131619599
$ws.Range("G6").Value = 0.8775636666666666
$ws.Range("H6").Value = 2.632691
$ws.Range("I6").Value = 0.1887436506618166
$ws.Range("J6").Value = 0.2083714858314108
$ws.Range("M6").Value = 0.7947223333333334
$ws.Range("N6").Value = 2.384167
$ws.Range("O6").Value = 0.05465238474931167
$ws.Range("P6").Value = 0.06446260031864033
$ws.Range("Q6").Value = 0.6974194448218889
$ws.Range("R6").Value = 6.276775003397
$ws.Range("S6").Value = 0.01031529061495928
$ws.Range("T6").Value = 0.01343216780895146
$ws.Range("G7").Value = 2.458038666666667
$ws.Range("H7").Value = 7.374116000000001
$ws.Range("I7").Value = 0.5286672739959656
$ws.Range("J7").Value = 0.5836444564186148
$ws.Range("M7").Value = 1.143899333333333
$ws.Range("N7").Value = 3.431698
$ws.Range("O7").Value = 0.0786649926114418
$ws.Range("P7").Value = 0.09278552072412605
$ws.Range("Q7").Value = 2.811748792107556
$ws.Range("R7").Value = 25.30573912896801
$ws.Range("S7").Value = 0.04158760720280371
$ws.Range("T7").Value = 0.05415375480655067
$ws.Range("G8").Value = 2.458038666666667
$ws.Range("H8").Value = 7.374116000000001
$ws.Range("I8").Value = 0.5286672739959656
$ws.Range("J8").Value = 0.5836444564186148
$ws.Range("O8").Value = 0.4082467010060103
$ws.Range("P8").Value = 0.4815278242489749
$ws.Range("S8").Value = 0.2158266705386935
$ws.Range("T8").Value = 0.2810410452342312
$ws.Range("G9").Value = 2.458038666666667
$ws.Range("H9").Value = 7.374116000000001
$ws.Range("I9").Value = 0.5286672739959656
$ws.Range("J9").Value = 0.5836444564186148
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.02736833333333333
$ws.Range("N9").Value = 0.082105
$ws.Range("O9").Value = 0.001882097206211744
$ws.Range("P9").Value = 0.002219937529192361
$ws.Range("Q9").Value = 0.06727242157555556
$ws.Range("R9").Value = 0.6054517941800001
$ws.Range("S9").Value = 0.0009950031994033852
$ws.Range("T9").Value = 0.001295654232508759
$ws.Range("G10").Value = 2.458038666666667
$ws.Range("H10").Value = 7.374116000000001
$ws.Range("I10").Value = 0.5286672739959656
$ws.Range("J10").Value = 0.5836444564186148
$ws.Range("M10").Value = 6.638933
$ws.Range("N10").Value = 13.277866
$ws.Range("O10").Value = 0.4565538244270245
$ws.Range("P10").Value = 0.3590041171790666
$ws.Range("Q10").Value = 16.31875401940933
$ws.Range("R10").Value = 97.91252411645601
$ws.Range("S10").Value = 0.2413650657922677
$ws.Range("T10").Value = 0.209530762823021
$ws.Range("G11").Value = 2.458038666666667
$ws.Range("H11").Value = 7.374116000000001
$ws.Range("I11").Value = 0.5286672739959656
$ws.Range("J11").Value = 0.5836444564186148
$ws.Range("M11").Value = 0.7947223333333334
$ws.Range("N11").Value = 2.384167
$ws.Range("O11").Value = 0.05465238474931167
$ws.Range("P11").Value = 0.06446260031864033
$ws.Range("Q11").Value = 1.953458224596889
$ws.Range("R11").Value = 17.581124021372
$ws.Range("S11").Value = 0.02889292726279729
$ws.Range("T11").Value = 0.03762323932230326
$ws.Range("G12").Value = 1.313898
$ws.Range("H12").Value = 2.627796
$ws.Range("I12").Value = 0.2825890753422177
$ws.Range("J12").Value = 0.2079840577499744
$ws.Range("M12").Value = 1.143899333333333
$ws.Range("N12").Value = 3.431698
$ws.Range("O12").Value = 0.0786649926114418
$ws.Range("P12").Value = 0.09278552072412605
$ws.Range("Q12").Value = 1.502967046268
$ws.Range("R12").Value = 9.017802277608
$ws.Range("S12").Value = 0.02222986752386973
$ws.Range("T12").Value = 0.01929790910064808
$ws.Range("G13").Value = 1.313898
$ws.Range("H13").Value = 2.627796
$ws.Range("I13").Value = 0.2825890753422177
$ws.Range("J13").Value = 0.2079840577499744
$ws.Range("O13").Value = 0.4082467010060103
$ws.Range("P13").Value = 0.4815278242489749
$ws.Range("Q13").Value = 7.799928761074001
$ws.Range("R13").Value = 46.79957256644401
$ws.Range("S13").Value = 0.1153660577487993
$ws.Range("T13").Value = 0.1001501108068183
$ws.Range("G14").Value = 1.313898
$ws.Range("H14").Value = 2.627796
$ws.Range("I14").Value = 0.2825890753422177
$ws.Range("J14").Value = 0.2079840577499744
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M14").Value = 0.02736833333333333
$ws.Range("N14").Value = 0.082105
$ws.Range("O14").Value = 0.001882097206211744
$ws.Range("P14").Value = 0.002219937529192361
$ws.Range("Q14").Value = 0.03595919843
$ws.Range("R14").Value = 0.21575519058
$ws.Range("S14").Value = 0.0005318601092075478
$ws.Range("T14").Value = 0.0004617116152728796
$ws.Range("G15").Value = 1.313898
$ws.Range("H15").Value = 2.627796
$ws.Range("I15").Value = 0.2825890753422177
$ws.Range("J15").Value = 0.2079840577499744
$ws.Range("M15").Value = 6.638933
$ws.Range("N15").Value = 13.277866
$ws.Range("O15").Value = 0.4565538244270245
$ws.Range("P15").Value = 0.3590041171790666
$ws.Range("Q15").Value = 8.722880790833999
$ws.Range("R15").Value = 34.891523163336
$ws.Range("S15").Value = 0.1290171230887861
$ws.Range("T15").Value = 0.07466713303984956
$ws.Range("G16").Value = 1.313898
$ws.Range("H16").Value = 2.627796
$ws.Range("I16").Value = 0.2825890753422177
$ws.Range("J16").Value = 0.2079840577499744
$ws.Range("M16").Value = 0.7947223333333334
$ws.Range("N16").Value = 2.384167
$ws.Range("O16").Value = 0.05465238474931167
$ws.Range("P16").Value = 0.06446260031864033
$ws.Range("Q16").Value = 1.044184084322
$ws.Range("R16").Value = 6.265104505932
$ws.Range("S16").Value = 0.01544416687155511
$ws.Range("T16").Value = 0.01340719318738561
